$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.779705762863159
$ws.Range("B1").Value = 3.573288679122925
$ws.Range("C1").Value = 2.695905447006226
$ws.Range("D1").Value = 2.560572147369385
$ws.Range("E1").Value = 2.600834369659424
